$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.276.64"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.811.02"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'312.64"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5132"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("D8").Value = "'0.3929"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").Value = "'0.07828"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'41.10"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'6.379"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'20.45"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "'7.344"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "1.807.08"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'92.82"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "'0.06584"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'6.016"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "28.325.23"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").Value = "'2.232"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "'160.01"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "'2.465"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").Value = "2.016.56"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "'127.62"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").Value = "'0.1097"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'1.060"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.582"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.652"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'0.07113"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "'9.150"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("D37").Value = "'0.02354"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "'0.2179"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").Value = "'11.63"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").Value = "'5.023"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "'0.6185"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'1.160"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("D44").Value = "'13.24"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "'0.5967"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "'1.306"
$ws.Range("E46").Value = "  -5.56%  "
$ws.Range("D47").Value = "'3.730"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "'125.38"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'1.212"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Value = "'0.06800"
